$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 8 (student #7)
$ws.Range("B8").Value = "Juan Diego Gonzalez Antoniazzi"
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = "jdgaprogrammer@gmail.com"
$ws.Range("E8").Value = "jdga1997"

# Add hyperlink for email cell D8
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:jdgaprogrammer@gmail.com")

# Match the style used by the other email hyperlink cells (D4:D7)
$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)

# Update selection to B10
$ws.Range("B10").Select()
